$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new task row above row 10 (shifts existing rows 10-31 down to 11-32),
# then remove one of the now-duplicated blank template rows (old row 15, now at
# row 16) so the sheet keeps 31 rows total.
$ws.Rows("10:10").Insert()
$ws.Rows("16:16").Delete()

# Fill in the new "Implementation" task row (row 10)
$ws.Range("B10").Value = "Implementation"
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3

# Update existing values
$ws.Range("H2").Value = 8
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 0.6
$ws.Range("F9").Value = 6
$ws.Range("G9").Value = 0.8

# Update the saved selection to match the author's final cursor position
$ws.Range("G11").Select()
